$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# --- Locate the ID placeholder run inside the first paragraph -------------
$idRange = $p1.Range.Duplicate
[void]$idRange.Find.Execute("**ID__AFFARS_5318_topic_5__ID**", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)

# Anything between the end of the placeholder text and the paragraph mark
# (e.g. the extra trailing " " run) is removed first, then the placeholder
# text itself is updated in place so its run/formatting is preserved.
$paraEnd = $p1.Range.End
$tailRange = $d.Range($idRange.End, $paraEnd - 1)
if ($tailRange.Start -lt $tailRange.End) {
    $tailRange.Delete()
}

$idRange.Text = "**ID__AFFARS_SUBPART_5318_2__ID**"

# --- Paragraph-level formatting changes on the first paragraph ------------
# Add a paragraph border with a 5pt space (offset) on all four sides.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25
